# Applies the "render website, remove theme (not needed) from docs"
# styles.xml changes to the document's style sheet:
#   1. Add a new "Abstract Title" paragraph style (based on Normal,
#      followed by Abstract).
#   2. Change the "Abstract" style's space-before from 15pt (300 twips)
#      to 5pt (100 twips).
#   3. Add a new "Footnote Block Text" paragraph style (based on
#      Footnote Text, followed by Footnote Text).

$d = $word.ActiveDocument

# --- 1. New "Abstract Title" style -----------------------------------
$abstractTitle = $d.Styles.Add("Abstract Title", 1)
$abstractTitle.BaseStyle = "Normal"
$abstractTitle.NextParagraphStyle = "Abstract"
$abstractTitle.QuickStyle = $true

$abstractTitle.ParagraphFormat.KeepWithNext = $true
$abstractTitle.ParagraphFormat.KeepTogether = $true
$abstractTitle.ParagraphFormat.Alignment = 1
$abstractTitle.ParagraphFormat.SpaceBefore = 15
$abstractTitle.ParagraphFormat.SpaceAfter = 0

$abstractTitle.Font.Size = 10
$abstractTitle.Font.SizeBi = 10
$abstractTitle.Font.Bold = $true
# RGB 34 5A 8A, packed as the 0x00BBGGRR long Word expects for Font.Color
$abstractTitle.Font.Color = 0x8A5A34

# --- 2. "Abstract" style: space-before 300 -> 100 twips (15pt -> 5pt) -
$abstract = $d.Styles("Abstract")
$abstract.ParagraphFormat.SpaceBefore = 5

# --- 3. New "Footnote Block Text" style -------------------------------
$footnoteBlockText = $d.Styles.Add("Footnote Block Text", 1)
$footnoteBlockText.BaseStyle = "Footnote Text"
$footnoteBlockText.NextParagraphStyle = "Footnote Text"
$footnoteBlockText.Priority = 9
$footnoteBlockText.UnhideWhenUsed = $true
$footnoteBlockText.QuickStyle = $true

$footnoteBlockText.ParagraphFormat.SpaceBefore = 5
$footnoteBlockText.ParagraphFormat.SpaceAfter = 5
$footnoteBlockText.ParagraphFormat.FirstLineIndent = 0
$footnoteBlockText.ParagraphFormat.LeftIndent = 24
$footnoteBlockText.ParagraphFormat.RightIndent = 24

Write-Output "Styles updated: Abstract Title (added), Abstract (spacing), Footnote Block Text (added)"
